$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "Modelo" in F1, reusing the same header style as A1:E1
# (bold, thin border, centered) by copying A1's format onto F1.
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "Modelo"

# Update the MSE / R2 / MAE values in row 2
$ws.Range("B2").Value = 0.2223578297307116
$ws.Range("C2").Value = 0.9969269545865439
$ws.Range("D2").Value = 0.3617619540350124

# Add the model name value in F2 (plain text, same cell formatting as E2)
$ws.Range("F2").Value = "Pipeline(steps=[('model', AdaBoostRegressor())])"
